# Update countries & provincias Spain
# Applies the daily COVID data refresh to the "Pais" sheet:
#   - bumps the "Datos actualizados..." timestamp in A1
#   - re-sorts Bolivia/Israel, Belgica/Kuwait/China and Jamaica/Estonia
#     (their totals changed enough to change rank order) together with
#     their refreshed case counts
#   - refreshes a few other countries' case counts in place

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $name, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $new
    $ws.Cells.Item($row, 4).Value = $active
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 04:59"

# Bolivia moves above Israel (row 30/31) with refreshed figures
Set-CountryRow 30 "Bolivia" 116598 630 60408 51163 0 61 5027
Set-CountryRow 31 "Israel" 116596 0 95009 20648 0 0 939

# Kazajistan (row 33) keeps its place, figures refreshed
Set-CountryRow 33 "Kazajistan" 105872 77 96297 8052 0 0 1523

# Belgica moves above Kuwait and China (rows 39/40/41) with refreshed figures
Set-CountryRow 39 "Belgica" 85236 194 18422 56919 0 1 9895
Set-CountryRow 40 "Kuwait" 85109 0 77224 7354 0 0 531
Set-CountryRow 41 "China" 85058 10 80208 216 0 0 4634

# Australia (row 72) keeps its place, figures refreshed
Set-CountryRow 72 "Australia" 25819 73 21350 3812 0 5 657

# Jamaica moves above Estonia (rows 136/137) with refreshed figures
Set-CountryRow 136 "Jamaica" 2459 102 890 1548 0 0 21
Set-CountryRow 137 "Estonia" 2375 0 2088 223 0 0 64
